$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'273.84"
$ws.Range("D3").Value = "'21.12"
$ws.Range("D4").Value = "'6.218"
$ws.Range("D6").Value = "'3.577"
$ws.Range("D7").Value = "'1.517"
$ws.Range("D8").Value = "'6.532"
$ws.Range("D10").Value = "'0.1650"
$ws.Range("D11").Value = "'0.08242"
$ws.Range("D12").Value = "'0.03413"
$ws.Range("D13").Value = "'0.03127"
$ws.Range("D14").Value = "'0.09128"
$ws.Range("D15").Value = "'3.767"
$ws.Range("D16").Value = "'0.001624"
$ws.Range("D17").Value = "'0.04673"
$ws.Range("D18").Value = "'0.006432"
$ws.Range("D19").Value = "'0.006129"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("D22").Value = "'3.720"
$ws.Range("D23").Value = "'2.318"
$ws.Range("D24").Value = "'0.01391"
$ws.Range("D25").Value = "'0.3323"
$ws.Range("D26").Value = "'0.1230"
$ws.Range("D28").Value = "'0.0002737"
$ws.Range("D40").Value = "'0.04759"
$ws.Range("B41").Value = "CEJI"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D41").Value = "'0.005502"
$ws.Range("E41").Value = "40CEJICEJIBestin24h"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.007036"
$ws.Range("E42").Value = "41KickTokenKICK"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1102"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.01027"
$ws.Range("D45").Value = "'0.00006090"
